$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 362, shifting existing rows 362:445 down to 363:446
$ws.Rows("362:362").Insert()

# Populate the newly inserted row 362 with a new price observation.
# Columns A,B,C,E,F,G,H,I,J,K,L,Q,R,T keep the same categorical values
# that the (now shifted-down) row used to have; D,M,N,O,P,S are new.
$ws.Cells.Item(362, 1).Value = 4
$ws.Cells.Item(362, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(362, 3).Value = "Los Lagos"
$ws.Cells.Item(362, 4).Value = 45135
$ws.Cells.Item(362, 5).Value = 10
$ws.Cells.Item(362, 6).Value = "Fruta"
$ws.Cells.Item(362, 7).Value = 100108
$ws.Cells.Item(362, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(362, 9).Value = 100108005
$ws.Cells.Item(362, 10).Value = "Piña"
$ws.Cells.Item(362, 11).Value = "Caramelo"
$ws.Cells.Item(362, 12).Value = "Segunda"
$ws.Cells.Item(362, 13).Value = 150
$ws.Cells.Item(362, 14).Value = 22000
$ws.Cells.Item(362, 15).Value = 22000
$ws.Cells.Item(362, 16).Value = 22000
$ws.Cells.Item(362, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(362, 18).Value = "Ecuador"
$ws.Cells.Item(362, 19).Value = 1571
$ws.Cells.Item(362, 20).Value = 14
